# The document was further edited near the "Assignment 8" title (code
# from previous programs was (re)loaded/typed in that area). Word's
# automatic "_GoBack" bookmark - which marks the location of the most
# recent edit - therefore moves from the end of the date line ("...,
# 2015") to right after "Assignment 8 " (just before "Writeup") on the
# title line.

$d = $word.ActiveDocument

# Drop the existing _GoBack bookmark (currently sitting at the end of
# the "November 2nd, 2015" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the title paragraph ("Assignment 8 Writeup") and split it right
# after "Assignment 8 " / right before "Writeup".
$titlePara = $d.Paragraphs(5)
$splitPoint = $titlePara.Range.Start + "Assignment 8 ".Length

$r = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $r)
